$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new trade row (row 8)
$ws.Range("A8").Value = 42649.644849537035
$ws.Range("A8").NumberFormat = "m/d/yy h:mm"

$ws.Range("B8").Value = $true

$ws.Range("C8").Value = 10068.83
$ws.Range("D8").Value = 9957.7999999999993
$ws.Range("E8").Value = 18.829999999999998
$ws.Range("F8").Value = 19.25

$ws.Range("G8").Value = $false
$ws.Range("G8").NumberFormat = "m/d/yy h:mm"

$ws.Range("H8").Value = 2.23

$ws.Range("I8").Value = $false

# Widen column C slightly to fit the new (wider) value in C8
$ws.Columns.Item(3).ColumnWidth = 8
